# ===========================================================================
# [ADDITIONAL SCRAPING] added scraping code for extra bowling attributes
# and excel sheets
#
# 1. Clean up stray empty placeholder cells on "ODI Batting Extra" that were
#    left over from the batting scrape (B/C/D/E columns with no data).
# 2. Add a new "ODI Bowling Extra" worksheet (mirrors "ODI Batting Extra")
#    holding MAIDEN_OVERS / PERCENT_WICKETS_OF_ALL per match.
# ===========================================================================

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: remove the stray blank cells from "ODI Batting Extra"
# ---------------------------------------------------------------------------
$battingExtra = $wb.Worksheets.Item("ODI Batting Extra")

$rangesToClear = @(
    "B2:E2",
    "B5:E5",
    "E6:E6",
    "E7:E7",
    "B13:E13",
    "B14:E14",
    "B15:E15",
    "B17:E17",
    "B21:E21"
)
foreach ($addr in $rangesToClear) {
    $battingExtra.Range($addr).ClearContents()
}

# ---------------------------------------------------------------------------
# Step 2: add the new "ODI Bowling Extra" worksheet after the last sheet
# ---------------------------------------------------------------------------
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$bowlingExtra = $wb.Worksheets.Add($null, $lastSheet)
$bowlingExtra.Name = "ODI Bowling Extra"

# Copy the header look & feel (bold / centered / bordered) from the sibling
# "ODI Batting Extra" sheet so the new sheet matches the workbook's style.
$battingExtra.Range("A1:C1").Copy()
$bowlingExtra.Range("A1:C1").PasteSpecial(-4122)  # xlPasteFormats
$bowlingExtra.Application.CutCopyMode = $false

$bowlingExtra.Range("A1").Value = "MATCH_CODE"
$bowlingExtra.Range("B1").Value = "MAIDEN_OVERS"
$bowlingExtra.Range("C1").Value = "PERCENT_WICKETS_OF_ALL"

# Data rows 2-21: MATCH_CODE, MAIDEN_OVERS, PERCENT_WICKETS_OF_ALL
$rows = @(
    @("3901", "0",    "20.00%"),
    @("3902", $null,  $null),
    @("3904", $null,  $null),
    @("3908", "0",    $null),
    @("3910", $null,  $null),
    @("3923", $null,  $null),
    @("3924", "0",    "20.00%"),
    @("3927", $null,  $null),
    @("4061", "0",    $null),
    @("4062", $null,  $null),
    @("4063", $null,  $null),
    @("4064", "2",    "20.00%"),
    @("4065", $null,  $null),
    @("4096", "2",    "10.00%"),
    @("4098", $null,  $null),
    @("4099", "0",    $null),
    @("4344", "0",    "10.00%"),
    @("4413", "0",    $null),
    @("4414", $null,  $null),
    @("4417", "0",    "40.00%")
)

$dataRange = $bowlingExtra.Range("A2:C21")
$dataRange.NumberFormat = "@"

$r = 2
foreach ($row in $rows) {
    $bowlingExtra.Cells.Item($r, 1).Value = $row[0]
    if ($null -ne $row[1]) {
        $bowlingExtra.Cells.Item($r, 2).Value = $row[1]
    }
    if ($null -ne $row[2]) {
        $bowlingExtra.Cells.Item($r, 3).Value = $row[2]
    }
    $r++
}

# Revert to the "no explicit style" look used by the rest of the data cells
# in this workbook (only the header row carries a named style).
$dataRange.Style = "Normal"

[void]$bowlingExtra.Range("A1").Select()
